# Update "CNN Data" session table: new Elapsed Time (s) / Enemy in Sight
# readings replace the previous (shorter) session, growing the used range
# from A1:B57 to A1:B80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$elapsedTime = @(
    3,4,5,6,7,8,9,10,10,11,
    12,13,14,15,16,17,17,18,19,20,
    21,22,23,23,24,25,26,27,28,29,
    30,30,31,32,33,34,35,36,37,37,
    38,39,40,41,42,43,44,44,45,46,
    47,48,49,50,51,51,52,53,54,55,
    56,57,58,58,59,60,61,62,63,64,
    65,65,66,67,68,69,70,71,71
)

$enemyInSight = @(
    $false,$false,$false,$false,$false,$false,$false,$false,$false,$true,
    $true,$false,$true,$false,$false,$false,$false,$false,$false,$true,
    $false,$true,$true,$false,$false,$false,$true,$false,$false,$false,
    $false,$false,$false,$true,$false,$false,$false,$false,$false,$false,
    $false,$false,$false,$false,$false,$true,$false,$false,$false,$false,
    $false,$false,$false,$false,$true,$false,$false,$false,$false,$true,
    $false,$false,$false,$false,$false,$true,$false,$true,$false,$false,
    $false,$false,$false,$false,$false,$false,$false,$false,$false
)

$rowCount = $elapsedTime.Length
$data = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $elapsedTime[$i]
    $data[$i,1] = $enemyInSight[$i]
}

$firstRow = 2
$lastRow = $firstRow + $rowCount - 1
$rangeAddress = "A" + $firstRow + ":B" + $lastRow
$ws.Range($rangeAddress).Value = $data
